# Update cryptos list - GitHub Actions scheduled refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the cell to be stored as text even when the value looks numeric
    # (e.g. "215.42"), without leaving a residual explicit style behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "26.023.61"

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.643.04"
$ws.Range("E3").Value = "  +0.27%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.27%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "215.42"
$ws.Range("E5").Value = "  +0.25%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.04%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.28%  "

# Row 8 - now Dogecoin (was Cardano)
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Range("D8") "0.0638"
$ws.Range("E8").Value = "  +0.12%  "

# Row 9 - now Cardano (was Dogecoin)
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue $ws.Range("D9") "0.255"
$ws.Range("E9").Value = "  +0.02%  "

# Row 10 - Solana
Set-TextValue $ws.Range("D10") "19.58"
$ws.Range("E10").Value = "  -0.40%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.09%  "

# Row 13 - WrappedEther
Set-TextValue $ws.Range("D13") "1.612.93"
$ws.Range("E13").Value = "  -1.29%  "

# Row 14 - Polygon
$ws.Range("E14").Value = "  -0.14%  "

# Row 15 - Litecoin
Set-TextValue $ws.Range("D15") "63.45"

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +0.21%  "

# Row 17 - WrappedBTC
Set-TextValue $ws.Range("D17") "26.058.16"
$ws.Range("E17").Value = "  +0.36%  "

# Row 19 - BitcoinCash
Set-TextValue $ws.Range("D19") "194.70"
$ws.Range("E19").Value = "  +0.25%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  -0.44%  "

# Row 21 - Avalanche
$ws.Range("E21").Value = "  -0.32%  "

# Row 22 - Chainlink
$ws.Range("E22").Value = "  -1.19%  "

# Row 23 - Stellar
$ws.Range("E23").Value = "  +4.63%  "

# Row 24 - Monero
Set-TextValue $ws.Range("D24") "143.98"
$ws.Range("E24").Value = "  -0.16%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  -0.56%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  +0.08%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  +0.58%  "

# Row 28 - EthereumClassic
Set-TextValue $ws.Range("D28") "15.50"
$ws.Range("E28").Value = "  +0.01%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  +0.26%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  -1.13%  "

# Row 31 - now Filecoin (was InternetComputer(DFINITY))
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D31") "3.26"
$ws.Range("E31").Value = "  +0.79%  "

# Row 32 - now InternetComputer(DFINITY) (was Filecoin)
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D32") "3.28"
$ws.Range("E32").Value = "  -0.76%  "

# Row 33 - LidoDAOToken
$ws.Range("E33").Value = "  -0.24%  "

# Row 35 - ARBITRUM
Set-TextValue $ws.Range("D35") "0.903"
$ws.Range("E35").Value = "  -0.08%  "

# Row 36 - Maker
Set-TextValue $ws.Range("D36") "1.130.29"

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  -1.46%  "

# Row 38 - MXToken
$ws.Range("E38").Value = "  -0.25%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  -0.26%  "

# Row 40 - now Quant (was FraxShare)
$ws.Range("B40").Value = "Quant"
$ws.Range("C40").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D40") "98.81"
$ws.Range("E40").Value = "  -0.45%  "

# Row 41 - now FraxShare (was Quant)
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D41") "5.43"
$ws.Range("E41").Value = "  +0.32%  "

# Row 42 - TrustWalletToken
Set-TextValue $ws.Range("D42") "0.794"
$ws.Range("E42").Value = "  -0.86%  "

# Row 43 - BabyDogeCoin
$ws.Range("E43").Value = "  +0.21%  "

# Row 44 - Aave
Set-TextValue $ws.Range("D44") "56.49"
$ws.Range("E44").Value = "  +0.02%  "

# Row 45 - RenderToken
$ws.Range("E45").Value = "  +2.82%  "

# Row 46 - Cronos
$ws.Range("E46").Value = "  -1.55%  "

# Row 47 - EnergySwap
$ws.Range("E47").Value = "  +2.11%  "

# Row 48 - Mantle
$ws.Range("E48").Value = "  -0.26%  "

# Row 49 - USDD
$ws.Range("E49").Value = "  +0.25%  "

# Row 50 - Algorand
Set-TextValue $ws.Range("D50") "0.0947"
$ws.Range("E50").Value = "  -1.63%  "

# Row 51 - Aptos
Set-TextValue $ws.Range("D51") "5.52"
$ws.Range("E51").Value = "  -0.18%  "
